$d = $word.ActiveDocument

function Split-At($doc, $pos, $endpos) {
    $r = $doc.Range($pos, $endpos)
    $r.Bold = 1
    $r.Bold = 0
}

# --- Change 1: add sentence to the end of the "Aktion verwalten" paragraph ---
$anchor1 = $d.Content
$anchor1.Find.Execute("Aktion löschen. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos1 = $anchor1.End
$ins1 = $d.Range($insPos1, $insPos1)
$ins1.InsertAfter("Es handelt sich um eine Verfeinerung.")
$splitEnd1 = $insPos1 + 37
Split-At $d $insPos1 $splitEnd1

# --- Change 2: insert new paragraphs after the "Letzlich..." paragraph ---
$lastp = $d.Paragraphs.Last
$lastRng = $lastp.Range
$lastRng.InsertParagraphAfter()
$cursor = $null

# paragraph group 1: style=Heading2
$newp = $d.Paragraphs.Last
$newp.Style = "Heading2"
$pStart = $newp.Range.Start
$newp.Range.Text = "Hilfsmittel verwalten"

# paragraph group 2: style=Normal
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart = $newp.Range.Start
$newp.Range.Text = "Der Use Case „Hilfsmittel“ ist ähnlich wie die vorrangehenden Use Cases ein Fall in dem eine Entität verwaltet wird. Das heißt es werden die Prozesse Anzeigen, Erstellen, Ändern und Löschen abgebildet."
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" Akteur ist in diesem Fall ein beliebiger Nutzer, der die entsprechende Berechtigung hat.")
$pEndNow = $newp.Range.End
$bpos1_0 = $pStart + 201
Split-At $d $bpos1_0 $pEndNow

# paragraph group 3: style=Normal
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart = $newp.Range.Start
$newp.Range.Text = "Das Diagramm besteht aus vier wesentlichen Bestandteilen: Hilfsmittel löschen, Hilfsmittel anzeigen, Hiflsmittel anlegen, sowie Hilfsmittel ändern."
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" Hilfsmittel anzeigen: Dieser Use Case zeigt eine Liste oder auch nur 1 Hilfsmittel an, dafür gibt es ein extend auf den Use Case „Liste anzeigen“.  Der Anwendungsfall Hilfsmittel löschen tritt auf wenn der Nutzer ein Hilfsmittel entfernen möchte, dafür wird ihm die Liste der Hilfsmittel angezeigt, weswegen der entsprechende Use Case ein include hat. Selbes Include gilt ebenfalls für den Anwendungsfall Hilfsmittel ändern, mit welchem der Nutzer die ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Attribute einer Entität Hilfsmittel verändern kann. Dafür wird ebenfalls eine Anzeige der Hilfsmittel benötigt, wo der korrespondierende Anwendungsfall wieder ins Spiel kommt. Weiterhin besitzt der Use Case includes auf folgende atomaren Anwendungsfälle, welche jeweils das entsprechende Attribut modifizieren: „Titel ändern“, „Beschreibung ändern“, „Kosten ändern“, „Art ändern“, sowie „Beleg ändern“. Letzterer hat ebenfalls Includes auf die Use Cases „Beleg hinzufügen“ und „Beleg entfernen“. Letztlich kann der Akteur ebenfalls ein Hilfsmittel anlegen, wofür der Anwendungsfall „Hilfsmittel anlegen“ existiert. Dieser hat wiederrum includes auf die atomaren Use Cases: „Titel eingeben“, „Beschreibung eingeben“, „Kosten eingeben“, „Art eingeben“ und „Beleg hinzufügen“")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(", welche jeweils das entsprechende Attribut der Entität setzen.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" ")
$pEndNow = $newp.Range.End
$bpos2_0 = $pStart + 1435
Split-At $d $bpos2_0 $pEndNow
$bpos2_1 = $pStart + 1372
Split-At $d $bpos2_1 $pEndNow
$bpos2_2 = $pStart + 600
Split-At $d $bpos2_2 $pEndNow
$bpos2_3 = $pStart + 147
Split-At $d $bpos2_3 $pEndNow

# paragraph group 4: style=Heading2
$newp = $d.Paragraphs.Last
$newp.Style = "Heading2"
$pStart = $newp.Range.Start
$newp.Range.Text = "Caterer verwalten"

# paragraph group 5: style=Normal
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart = $newp.Range.Start
$newp.Range.Text = "Das Diagramm Caterer verwalten bildet den Use Case „Caterer verwalten“  ab"
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(". ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Die Verwaltung der Attribute Name sowie Beschreibung wurden zur Wahrung der Übersichtlichkeit vernachlässigt.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Akteur ist in diesem Fall ein beliebiger Nutzer, der die entsprechende Berechtigung hat.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Es werden allen möglichen Anwendungsfälle dargestellt, die im Rahmen der Verwaltung von der Entität „Caterer“ auftreten können. Dazu gehören maßgeblich die Fälle: „Caterer anzeigen“, „Caterer hinzufügen“, „Caterer ändern“, sowie Caterer löschen“.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" Der Anwendungsfall Caterer anzeigen hat eine extend auf Liste anzeigen, da er eine List von Catereren anzeigt. Dafür benötigt er Includes auf Use Case „Anzeigen“ der entsprechenden Enitätsmengen, die er als Attribute besitzt. Dazu gehören: „Beleg anzeigen“, „Essen anzeigen“, „Trinken anzeigen“. Der Anwendungsfall „Caterer löschen“ includiert „Caterer anzeigen“, da zum Löschen der entsprechenden Entität sie zunächst dem Nutzer angezeigt werden muss. Ein weiter Use Case ist Catere hinzufügen, in welchem der Akteur eine neue  Entität vom Typ „Caterer“ erstellt. Dafür besitzt dieser(Use Case) includes auf „Beleg hinzufügen“, „Essen hinzufügen“ und „Trinken hinzufügen“. Diese Use Cases wiederum besitzen ein include auf die jeweils korrespondierenden „Anzeigen“ Anwendungsfall, da die entsprechenden Entität vor dem hinzufügen angezeigt werden müssen. Der letzte Use Case in diesem Diagramm ist „Caterer ändern“ welcher ein extends auf Caterer hinzufügen hat, da er diesen erweitert. Er erhält zusätzliche includes auf „Beleg entfernen“, „Essen entfernen“ sowie „Trinken entfernen“ , diese hätten jeweils auch ein Include auf den entsprechenden „Anzeigen“ Use Case, diese wurden aber aufgrund der Übersichtlichkeit nicht eingefügt.")
$pEndNow = $newp.Range.End
$bpos4_0 = $pStart + 521
Split-At $d $bpos4_0 $pEndNow
$bpos4_1 = $pStart + 275
Split-At $d $bpos4_1 $pEndNow
$bpos4_2 = $pStart + 274
Split-At $d $bpos4_2 $pEndNow
$bpos4_3 = $pStart + 186
Split-At $d $bpos4_3 $pEndNow
$bpos4_4 = $pStart + 185
Split-At $d $bpos4_4 $pEndNow
$bpos4_5 = $pStart + 76
Split-At $d $bpos4_5 $pEndNow
$bpos4_6 = $pStart + 74
Split-At $d $bpos4_6 $pEndNow

# paragraph group 6: style=Heading2
$newp = $d.Paragraphs.Last
$newp.Style = "Heading2"
$pStart = $newp.Range.Start
$newp.Range.Text = "Lebensmittel verwalten"

# paragraph group 7: style=Normal
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart = $newp.Range.Start
$newp.Range.Text = "Das D"
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("iagramm „Lebensmittel verwalten“")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" bildet den Use Case „Lebensmittel verwalten“  ab.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Akteur ist in diesem Fall ein beliebiger Nutzer, der die entsprechende Berechtigung hat. Es werden allen möglichen Anwendungsfälle dargestellt, die im Rahmen der Verwaltung von der Entität „")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Lebensmittel")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("“ auftreten können. Dazu gehören maßgeblich die Fälle: „")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Lebensmittel")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" anzeigen“, „")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Lebensmittel")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" hinzufügen“, „")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Lebensmittel")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" ändern“")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(", sowie „Lebensmittel")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" löschen“.")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter(" Der Use Case „Lebensmittel anzeigen“ hat ein extend auf List anzeigen, da es eine Liste von Lebensmitteln anzeigen kann. Diesen includiert der Anwendungsfall „Lebensmittel entfernen“, da vor dem Löschen einer Entität diese zunächst angezeigt werden muss.  Weiterhin existiert der Use Case Lebensmittel hinzufügen, welcher den Prozess abbildet wenn der Nutzer eine neue Entität vom Typ „Lebensmittel“ erstellen möchte. Dafür hat dieser Anwendungsfall includes auf díe atomaren Fälle „Titel hinzufügen“, „Beschreibung hinzufügen“, „Menge angeben“ sowie „Mengenbeschreibung angeben“. Diese setzten jeweils das entsprechende Attribut der Entität. Ähnlich ist der Use Case „Lebensmittel ändern“ aufgebaut. ")
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("Jedoch erstellt dieser kein neues Objekt, sondern modifiziert ein bestehendes. Dafür besitzt er die entsprechenden Includes auf die Andwendungsfälle „Titel ändern“, „Beschreibung ändern“, „Menge ändern“ sowie „Mengenbeschreibung ändern“. Diese sind ändern jeweils das korrespondiernde Attribut der Entität.")
$bmPos = $d.Content
$bmPos.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmPos)
$e = $d.Content
$e.Collapse(0)
$e.InsertAfter("  ")
$pEndNow = $newp.Range.End
$bpos6_0 = $pStart + 1457
Split-At $d $bpos6_0 $pEndNow
$bpos6_1 = $pStart + 1151
Split-At $d $bpos6_1 $pEndNow
$bpos6_2 = $pStart + 449
Split-At $d $bpos6_2 $pEndNow
$bpos6_3 = $pStart + 439
Split-At $d $bpos6_3 $pEndNow
$bpos6_4 = $pStart + 418
Split-At $d $bpos6_4 $pEndNow
$bpos6_5 = $pStart + 410
Split-At $d $bpos6_5 $pEndNow
$bpos6_6 = $pStart + 398
Split-At $d $bpos6_6 $pEndNow
$bpos6_7 = $pStart + 383
Split-At $d $bpos6_7 $pEndNow
$bpos6_8 = $pStart + 371
Split-At $d $bpos6_8 $pEndNow
$bpos6_9 = $pStart + 358
Split-At $d $bpos6_9 $pEndNow
$bpos6_10 = $pStart + 346
Split-At $d $bpos6_10 $pEndNow
$bpos6_11 = $pStart + 290
Split-At $d $bpos6_11 $pEndNow
$bpos6_12 = $pStart + 278
Split-At $d $bpos6_12 $pEndNow
$bpos6_13 = $pStart + 88
Split-At $d $bpos6_13 $pEndNow
$bpos6_14 = $pStart + 87
Split-At $d $bpos6_14 $pEndNow
$bpos6_15 = $pStart + 37
Split-At $d $bpos6_15 $pEndNow
$bpos6_16 = $pStart + 5
Split-At $d $bpos6_16 $pEndNow
